$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = ""
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 36
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.0'
$ws.Range("E8").Value = 'Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it''s ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F8").Value = 0
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '0.00'
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = ""

# Row 9
$ws.Range("A9").Value = 'P. point'
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 37
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3'
$ws.Range("E9").Value = 'Medium point (up to 6 mtr.)'
$ws.Range("F9").Value = 472
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '17464.00'
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = ""

# Row 10
$ws.Range("A10").Value = 'P. point'
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 52
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '4'
$ws.Range("E10").Value = 'Long point  (up to 10 mtr.)'
$ws.Range("F10").Value = 662
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '34424.00'
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = ""

# Row 11
$ws.Range("A11").Value = ""
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 29
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '2.0'
$ws.Range("E11").Value = 'Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it''s  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F11").Value = 0
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '0.00'
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = ""

# Row 12
$ws.Range("A12").Value = 'P. point'
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 78
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6'
$ws.Range("E12").Value = 'On board'
$ws.Range("F12").Value = 136
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '10608.00'
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = ""

# Row 13
$ws.Range("A13").Value = 'Each'
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 71
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.0'
$ws.Range("E13").Value = 'P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F13").Value = 23
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '1633.00'
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = ""

# Row 14
$ws.Range("A14").Value = 'Each'
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 72
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.0'
$ws.Range("E14").Value = 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F14").Value = 50
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '3600.00'
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = ""

# Row 15
$ws.Range("A15").Value = 'Each'
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 50
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.0'
$ws.Range("E15").Value = 'Providing & Fixing of  of 3/5 pin 6 amp. flush type  non modular socket  made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F15").Value = 33
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '1650.00'
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = ""

# Row 16
$ws.Range("A16").Value = 'Each'
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 28
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.0'
$ws.Range("E16").Value = 'Providing & Fixing of  ISI marked (IS:371) 6 amp surface type 3 pin ceiling rose with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screws including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F16").Value = 30
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '840.00'
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = ""

# Row 17
$ws.Range("A17").Value = 'Each'
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 67
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.0'
$ws.Range("E17").Value = 'Providing & Fixing of ISI marked (IS:1258) batten/angle lamp  holder with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screwsincluding making connection testing etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F17").Value = 30
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '2010.00'
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = ""

# Row 18
$ws.Range("A18").Value = 'Each'
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 32
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.0'
$ws.Range("E18").Value = 'Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F18").Value = 219
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '7008.00'
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = ""

# Row 19
$ws.Range("A19").Value = 'Each'
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 70
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.0'
$ws.Range("E19").Value = 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F19").Value = 303
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '21210.00'
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = ""

# Row 20
$ws.Range("A20").Value = ""
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 82
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.0'
$ws.Range("E20").Value = 'S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F20").Value = 0
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '0.00'
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = ""

# Row 21
$ws.Range("A21").Value = 'R. mtr.'
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 42
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16'
$ws.Range("E21").Value = '20 mm'
$ws.Range("F21").Value = 40
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '1680.00'
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = ""

# Row 22
$ws.Range("A22").Value = 'R. mtr.'
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 76
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17'
$ws.Range("E22").Value = '25 mm'
$ws.Range("F22").Value = 56
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '4256.00'
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = ""

# Row 23
$ws.Range("A23").Value = ""
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 49
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.0'
$ws.Range("E23").Value = 'Supplying and drawing FR PVC insulated & unsheathed flexible copper conductor as per PWD specification for electrical Works with ISI marked (IS:694) and as per IS 8130 : 2013 of 1.1 kV grade . Wire should be made from  99.90 % purity copper, class 2 stranding in acc. to IS:8130/IEC 60228 for  lower watt loss , oxygen free for less chances of oxidization, insulation PVC type A/C/D , flame retardant as per IS 10810-53, better amperage rating as per IS:3961 part 5,  in existing  surface or recessed PVC/ MS conduit/casing capping making connections with Copper Lugs of suitable size, Ferrules,testing etc. as required. OEM Must have its own in house NABL lab setup for all testing facilities for wires.   For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F23").Value = 0
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '0.00'
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = ""

# Row 24
$ws.Range("A24").Value = 'Mtr.'
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = 3
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '19'
$ws.Range("E24").Value = '2 x 2.5 sq. mm. + 1x1.5sqmm'
$ws.Range("F24").Value = 81
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '243.00'
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = ""

# Row 25
$ws.Range("A25").Value = 'Mtr.'
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = 66
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '20'
$ws.Range("E25").Value = '2 x 4.0 sq. mm. + 1 x 2.5 sq. mm.'
$ws.Range("F25").Value = 122
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '8052.00'
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = ""

# Row 26
$ws.Range("A26").Value = 'Set'
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 20
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.0'
$ws.Range("E26").Value = 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F26").Value = 5733
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '114660.00'
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = ""

# Row 27
$ws.Range("A27").Value = ""
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 57
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.0'
$ws.Range("E27").Value = 'Supply & Laying following size earth wire in horizontal or vertical run in ground/surface/recess including riveting, soldering, saddles,  making connection with GI/Cu purity purity >95%  thimble etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F27").Value = 0
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '0.00'
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = ""

# Row 28
$ws.Range("A28").Value = 'Mtr.'
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 55
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23'
$ws.Range("E28").Value = '8 SWG G.I. ( Hot Dipped  ) Wire '
$ws.Range("F28").Value = 20
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '1100.00'
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = ""

# Row 29
$ws.Range("A29").Value = ""
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 64
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.0'
$ws.Range("E29").Value = 'Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F29").Value = 0
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '0.00'
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = ""

# Row 30
$ws.Range("A30").Value = 'Each'
$ws.Range("B30").Value = 0
$ws.Range("C30").Value = 72
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '25'
$ws.Range("E30").Value = '1200 mm Sweep BEE 1 Star rated (service value >=4.0 to < 4.5 )'
$ws.Range("F30").Value = 1890
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '136080.00'
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = ""

# Row 31
$ws.Range("A31").Value = ""
$ws.Range("B31").Value = 0
$ws.Range("C31").Value = 61
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '16.0'
$ws.Range("E31").Value = 'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F31").Value = 0
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '0.00'
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = ""

# Row 32
$ws.Range("A32").Value = 'Each'
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 55
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27'
$ws.Range("E32").Value = '1170mm(+/-10%) LED batten with min. lumen output 2200 lm'
$ws.Range("F32").Value = 492
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '27060.00'
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = ""

# Row 33
$ws.Range("A33").Value = ""
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = 97
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '17.0'
$ws.Range("E33").Value = 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F33").Value = 0
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '0.00'
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = ""

# Row 34
$ws.Range("A34").Value = ""
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = 33
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29'
$ws.Range("E34").Value = 'Single pole MCB   (With B/C curve tripping Characteristics)'
$ws.Range("F34").Value = 0
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '0.00'
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = ""

# Row 35
$ws.Range("A35").Value = 'Each'
$ws.Range("B35").Value = 0
$ws.Range("C35").Value = 6
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '30'
$ws.Range("E35").Value = ' 6 A to 32 A rating'
$ws.Range("F35").Value = 187
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '1122.00'
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = ""

# Row 36
$ws.Range("A36").Value = ""
$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 17
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31'
$ws.Range("E36").Value = 'Double pole MCB(With B/C curve tripping Characteristics)'
$ws.Range("F36").Value = 0
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '0.00'
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = ""

# Row 37
$ws.Range("A37").Value = 'Each'
$ws.Range("B37").Value = 0
$ws.Range("C37").Value = 96
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '32'
$ws.Range("E37").Value = ' 50/63 A rating'
$ws.Range("F37").Value = 900
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '86400.00'
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = ""

# Row 38
$ws.Range("A38").Value = ""
$ws.Range("B38").Value = 0
$ws.Range("C38").Value = 84
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.0'
$ws.Range("E38").Value = 'Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F38").Value = 0
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '0.00'
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = ""

# Row 39
$ws.Range("A39").Value = ""
$ws.Range("B39").Value = 0
$ws.Range("C39").Value = 6
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34'
$ws.Range("E39").Value = 'Metal door (single phase) IK-09 and IP-43 with Metal end box'
$ws.Range("F39").Value = 0
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '0.00'
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = ""

# Row 40
$ws.Range("A40").Value = 'Each'
$ws.Range("B40").Value = 0
$ws.Range("C40").Value = 22
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '35'
$ws.Range("E40").Value = '8 Way (8+2)'
$ws.Range("F40").Value = 2184
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '48048.00'
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = ""

# Row 41
$ws.Range("A41").Value = ""
$ws.Range("B41").Value = 0
$ws.Range("C41").Value = 3
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36'
$ws.Range("E41").Value = 'Total'
$ws.Range("F41").Value = 0
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '0.00'
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = ""

# Row 42
$ws.Range("A42").Value = '%'
$ws.Range("B42").Value = 0
$ws.Range("C42").Value = 4
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '37'
$ws.Range("E42").Value = 'Add Tender Premium '
$ws.Range("F42").Value = 0
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '0.00'
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = ""

# Row 43
$ws.Range("A43").Value = ""
$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 94
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '38'
$ws.Range("E43").Value = 'Grand Total'
$ws.Range("F43").Value = 0
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '0.00'
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = ""

# Row 44
$ws.Range("A44").Value = ""

# Row 45
$ws.Range("A45").Value = ""
$ws.Range("B45").Value = ""
$ws.Range("C45").Value = ""
$ws.Range("D45").Value = ""
$ws.Range("E45").Value = 'Grand Total Rs.'
$ws.Range("F45").Value = ""
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '529148.00'
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = '529148.00'
$ws.Range("I45").Value = ""

# Row 46
$ws.Range("A46").Value = ""
$ws.Range("B46").Value = ""
$ws.Range("C46").Value = ""
$ws.Range("D46").Value = ""
$ws.Range("E46").Value = 'Tender Premium @ 0%'
$ws.Range("F46").Value = ""
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '0.00'
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = '0.00'
$ws.Range("I46").Value = ""

# Row 47
$ws.Range("A47").Value = ""
$ws.Range("B47").Value = ""
$ws.Range("C47").Value = ""
$ws.Range("D47").Value = ""
$ws.Range("E47").Value = 'NET PAYABLE AMOUNT Rs.'
$ws.Range("F47").Value = ""
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '529148.00'
$ws.Range("H47").NumberFormat = "@"
$ws.Range("H47").Value = '529148.00'
$ws.Range("I47").Value = ""
